$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H88").Value = 3975.5833
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3975.5833
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3975.5833
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = -4787.5833
$ws.Range("H91").Value = 3975.5833
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3975.5833
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3975.5833
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = -6783.5833
$ws.Range("H98").Value = 746.1111
$ws.Range("I98").Value = 746.1111
$ws.Range("K98").Value = 746.1111
$ws.Range("M98").Value = 751.8889
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null
$ws.Range("H122").Value = 746.1111
$ws.Range("I122").Value = 746.1111
$ws.Range("K122").Value = 2238.3333
$ws.Range("M122").Value = 211.6667000000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21140.28
$ws.Range("I32").Value = 3637.6843
$ws.Range("J32").Value = 154160
$ws.Range("K32").Value = 3637.6843
$ws.Range("L32").Value = 154160
$ws.Range("M32").Value = -3350.6843
$ws.Range("N32").Value = -154734
$ws.Range("H98").Value = 5649.75
$ws.Range("J98").Value = 5649.75
$ws.Range("L98").Value = 5649.75
$ws.Range("N98").Value = -11639.75
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").Value = $null
$ws.Range("H122").Value = 2471.35
$ws.Range("I122").Value = 2628.2222
$ws.Range("K122").Value = 7884.6666
$ws.Range("M122").Value = -5434.6666
$ws.Range("H132").Value = 2717.9
$ws.Range("I132").Value = 2592.2856
$ws.Range("J132").Value = 3597.2
$ws.Range("K132").Value = 7776.8568
$ws.Range("L132").Value = 10791.6
$ws.Range("M132").Value = -5246.8568
$ws.Range("N132").Value = -15851.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66245.82000000001
$ws.Range("I86").Value = 79880
$ws.Range("J86").Value = 2619.6667
$ws.Range("K86").Value = 79880
$ws.Range("L86").Value = 2619.6667
$ws.Range("M86").Value = -78757
$ws.Range("N86").Value = -4865.6667
$ws.Range("H89").Value = 66245.82000000001
$ws.Range("I89").Value = 79880
$ws.Range("J89").Value = 2619.6667
$ws.Range("K89").Value = 399400
$ws.Range("L89").Value = 13098.3335
$ws.Range("M89").Value = -393784
$ws.Range("N89").Value = -24330.3335
$ws.Range("H105").Value = 101664.95
$ws.Range("I105").Value = 92340.91
$ws.Range("J105").Value = 113061
$ws.Range("K105").Value = 92340.91
$ws.Range("L105").Value = 113061
$ws.Range("M105").Value = -90593.91
$ws.Range("N105").Value = -116555

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23579.762
$ws.Range("I31").Value = 26891.514
$ws.Range("J31").Value = 5128.5713
$ws.Range("K31").Value = 26891.514
$ws.Range("L31").Value = 5128.5713
$ws.Range("M31").Value = -26596.514
$ws.Range("N31").Value = -5718.5713
$ws.Range("H34").Value = 23579.762
$ws.Range("I34").Value = 26891.514
$ws.Range("J34").Value = 5128.5713
$ws.Range("K34").Value = 26891.514
$ws.Range("L34").Value = 5128.5713
$ws.Range("M34").Value = -26689.514
$ws.Range("N34").Value = -5532.5713
$ws.Range("H45").Value = 13750
$ws.Range("J45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -16186
$ws.Range("H58").Value = 7319.2
$ws.Range("I58").Value = 1134.9656
$ws.Range("J58").Value = 23623.092
$ws.Range("K58").Value = 1134.9656
$ws.Range("L58").Value = 23623.092
$ws.Range("M58").Value = -931.9656
$ws.Range("N58").Value = -24029.092
$ws.Range("H136").Value = 7319.2
$ws.Range("I136").Value = 1134.9656
$ws.Range("J136").Value = 23623.092
$ws.Range("K136").Value = 3404.8968
$ws.Range("L136").Value = 70869.276
$ws.Range("M136").Value = -854.8968
$ws.Range("N136").Value = -75969.276

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 100001110
$ws.Range("J25").Value = 111112180
$ws.Range("L25").Value = 333336540
$ws.Range("N25").Value = -333336878
$ws.Range("H30").Value = 100001110
$ws.Range("J30").Value = 111112180
$ws.Range("L30").Value = 333336540
$ws.Range("N30").Value = -333336744
$ws.Range("H131").Value = 817.15
$ws.Range("I131").Value = 454.53845
$ws.Range("J131").Value = 871.3333
$ws.Range("K131").Value = 1363.61535
$ws.Range("L131").Value = 2613.9999
$ws.Range("M131").Value = 3676.38465
$ws.Range("N131").Value = -12693.9999
$ws.Range("H139").Value = 2443.8462
$ws.Range("I139").Value = 1488
$ws.Range("J139").Value = 3041.25
$ws.Range("K139").Value = 4464
$ws.Range("L139").Value = 9123.75
$ws.Range("M139").Value = 676
$ws.Range("N139").Value = -19403.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 54645.676
$ws.Range("I70").Value = 91136.87
$ws.Range("J70").Value = 5275.2354
$ws.Range("K70").Value = 91136.87
$ws.Range("L70").Value = 5275.2354
$ws.Range("M70").Value = -90866.87
$ws.Range("N70").Value = -5815.2354
$ws.Range("H73").Value = 54645.676
$ws.Range("I73").Value = 91136.87
$ws.Range("J73").Value = 5275.2354
$ws.Range("K73").Value = 91136.87
$ws.Range("L73").Value = 5275.2354
$ws.Range("M73").Value = -90200.87
$ws.Range("N73").Value = -7147.2354
$ws.Range("H104").Value = 45550
$ws.Range("J104").Value = 45550
$ws.Range("L104").Value = 45550
$ws.Range("N104").Value = -52538
$ws.Range("H132").Value = 2346.3333
$ws.Range("I132").Value = 1723.5294
$ws.Range("J132").Value = 3858.8572
$ws.Range("K132").Value = 5170.5882
$ws.Range("L132").Value = 11576.5716
$ws.Range("M132").Value = -2640.5882
$ws.Range("N132").Value = -16636.5716

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 64350
$ws.Range("I40").Value = 201020.8
$ws.Range("K40").Value = 201020.8
$ws.Range("M40").Value = -200884.8
$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2219.6316
$ws.Range("J126").Value = 1574.75
$ws.Range("L126").Value = 4724.25
$ws.Range("N126").Value = -9664.25
$ws.Range("H132").Value = 5015.231
$ws.Range("I132").Value = 5240
$ws.Range("J132").Value = 4266
$ws.Range("K132").Value = 15720
$ws.Range("L132").Value = 12798
$ws.Range("M132").Value = -13190
$ws.Range("N132").Value = -17858
$ws.Range("H136").Value = 1311.1111
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1700
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 5100
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -10200
